$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 4047.9048
$ws.Range("I51").Value = 1727.6364
$ws.Range("K51").Value = 1727.6364
$ws.Range("M51").Value = -1243.6364
$ws.Range("H98").Value = 2952.3635
$ws.Range("I98").Value = 2388.5715
$ws.Range("J98").Value = 3939
$ws.Range("K98").Value = 2388.5715
$ws.Range("L98").Value = 3939
$ws.Range("M98").Value = -890.5715
$ws.Range("N98").Value = -6935
$ws.Range("H122").Value = 2952.3635
$ws.Range("I122").Value = 2388.5715
$ws.Range("J122").Value = 3939
$ws.Range("K122").Value = 7165.7145
$ws.Range("L122").Value = 11817
$ws.Range("M122").Value = -4715.7145
$ws.Range("N122").Value = -16717
$ws.Range("H137").Value = 2527530.2
$ws.Range("I137").Value = 2006.3529
$ws.Range("K137").Value = 6019.0587
$ws.Range("M137").Value = -3469.0587

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1638.4286
$ws.Range("I61").Value = 1575.8334
$ws.Range("J61").Value = 2014
$ws.Range("K61").Value = 1575.8334
$ws.Range("L61").Value = 2014
$ws.Range("M61").Value = -1363.8334
$ws.Range("H74").Value = 62826.777
$ws.Range("I74").Value = 84892.836
$ws.Range("J74").Value = 18694.666
$ws.Range("K74").Value = 84892.836
$ws.Range("L74").Value = 18694.666
$ws.Range("M74").Value = -84018.836
$ws.Range("N74").Value = -20442.666
$ws.Range("H77").Value = 62826.777
$ws.Range("I77").Value = 84892.836
$ws.Range("J77").Value = 18694.666
$ws.Range("K77").Value = 424464.18
$ws.Range("L77").Value = 93473.33
$ws.Range("M77").Value = -420096.18
$ws.Range("N77").Value = -102209.33
$ws.Range("H132").Value = 2866.5908
$ws.Range("I132").Value = 2127.4375
$ws.Range("J132").Value = 4837.6665
$ws.Range("K132").Value = 6382.3125
$ws.Range("L132").Value = 14512.9995
$ws.Range("M132").Value = -3852.3125
$ws.Range("N132").Value = -19572.9995
$ws.Range("H136").Value = 1638.4286
$ws.Range("I136").Value = 1575.8334
$ws.Range("J136").Value = 2014
$ws.Range("K136").Value = 4727.5002
$ws.Range("L136").Value = 6042
$ws.Range("M136").Value = -2177.5002

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 2139.5625
$ws.Range("I99").Value = 1777.75
$ws.Range("J99").Value = 2260.1667
$ws.Range("K99").Value = 1777.75
$ws.Range("L99").Value = 2260.1667
$ws.Range("M99").Value = -279.75
$ws.Range("N99").Value = -5256.1667
$ws.Range("H119").Value = 29868.75
$ws.Range("J119").Value = 29868.75
$ws.Range("L119").Value = 29868.75
$ws.Range("N119").Value = -39544.75
$ws.Range("H134").Value = 4583.485
$ws.Range("I134").Value = 4388.174
$ws.Range("J134").Value = 5032.7
$ws.Range("K134").Value = 13164.522
$ws.Range("L134").Value = 15098.1
$ws.Range("M134").Value = -10629.522
$ws.Range("N134").Value = -20168.1

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1527.1111
$ws.Range("I31").Value = 1468
$ws.Range("J31").Value = 2000
$ws.Range("K31").Value = 1468
$ws.Range("L31").Value = 2000
$ws.Range("M31").Value = -1173
$ws.Range("H34").Value = 1527.1111
$ws.Range("I34").Value = 1468
$ws.Range("J34").Value = 2000
$ws.Range("K34").Value = 1468
$ws.Range("L34").Value = 2000
$ws.Range("M34").Value = -1266
$ws.Range("H134").Value = 2944.4707
$ws.Range("I134").Value = 2918.6428
$ws.Range("K134").Value = 8755.928400000001
$ws.Range("M134").Value = -6220.928400000001
$ws.Range("H140").Value = 39540
$ws.Range("J140").Value = 39540
$ws.Range("L140").Value = 39540

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 653.7826
$ws.Range("I113").Value = 613.6
$ws.Range("J113").Value = 729.125
$ws.Range("K113").Value = 1840.8
$ws.Range("L113").Value = 2187.375
$ws.Range("M113").Value = 329.1999999999998
$ws.Range("N113").Value = -6527.375
$ws.Range("H131").Value = 884.63
$ws.Range("I131").Value = 550
$ws.Range("J131").Value = 905.9894
$ws.Range("K131").Value = 1650
$ws.Range("L131").Value = 2717.9682
$ws.Range("N131").Value = -12797.9682
$ws.Range("H132").Value = 2846.261
$ws.Range("I132").Value = 2403.65
$ws.Range("K132").Value = 21632.85
$ws.Range("M132").Value = -19102.85

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3410.5557
$ws.Range("I132").Value = 3563
$ws.Range("K132").Value = 10689
$ws.Range("M132").Value = -8159

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1878.7142
$ws.Range("I7").Value = 1754.7273
$ws.Range("J7").Value = 2333.3333
$ws.Range("K7").Value = 1754.7273
$ws.Range("L7").Value = 2333.3333
$ws.Range("M7").Value = -1642.7273
$ws.Range("N7").Value = -2557.3333
$ws.Range("H40").Value = 3324.75
$ws.Range("I40").Value = 3324.75
$ws.Range("K40").Value = 3324.75
$ws.Range("M40").Value = -3188.75
$ws.Range("H93").Value = 38202.21
$ws.Range("I93").Value = 1495.9375
$ws.Range("J93").Value = 233969
$ws.Range("K93").Value = 1495.9375
$ws.Range("L93").Value = 233969
$ws.Range("M93").Value = -247.9375
$ws.Range("N93").Value = -236465
$ws.Range("H126").Value = 1878.7142
$ws.Range("I126").Value = 1754.7273
$ws.Range("J126").Value = 2333.3333
$ws.Range("K126").Value = 5264.1819
$ws.Range("L126").Value = 6999.999899999999
$ws.Range("M126").Value = -2794.1819
$ws.Range("N126").Value = -11939.9999
$ws.Range("H132").Value = 3974.0386
$ws.Range("I132").Value = 4073
$ws.Range("K132").Value = 12219
$ws.Range("M132").Value = -9689

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3927.0667
$ws.Range("I132").Value = 4373.8945
$ws.Range("J132").Value = 3155.2727
$ws.Range("K132").Value = 13121.6835
$ws.Range("L132").Value = 9465.8181
$ws.Range("M132").Value = -10591.6835
$ws.Range("N132").Value = -14525.8181
$ws.Range("H135").Value = 34250
$ws.Range("J135").Value = 34250
$ws.Range("L135").Value = 34250
$ws.Range("N135").Value = -44390
$ws.Range("H136").Value = 2864.8857
$ws.Range("I136").Value = 2866.6428
$ws.Range("K136").Value = 8599.928400000001
$ws.Range("M136").Value = -6049.928400000001
